{"js": "// Update the worksheet date and the 25 multiplication problems.\n// Row 0 (paragraph) holds the date string; the table's rows 0, 4, 9, 14, 19\n// (0-indexed) are the only rows that carry the \"AxB=\" text cells \u2014 the rows\n// in between are blank spacer rows.\n\n// 1) Update the date/day-of-week heading.\nconst body = context.document.body;\nconst dateResults = body.search(\"2023-11-10 Friday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2023-11-11 Saturday\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Update the multiplication problems table.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\nconst rowUpdates = {\n  0: [\"46\u00d711=\", \"93\u00d745=\", \"27\u00d782=\", \"35\u00d775=\", \"87\u00d763=\"],\n  4: [\"19\u00d778=\", \"26\u00d769=\", \"17\u00d784=\", \"12\u00d744=\", \"71\u00d762=\"],\n  9: [\"45\u00d732=\", \"74\u00d713=\", \"88\u00d772=\", \"29\u00d757=\", \"18\u00d752=\"],\n  14: [\"97\u00d750=\", \"96\u00d798=\", \"25\u00d741=\", \"81\u00d744=\", \"75\u00d779=\"],\n  19: [\"37\u00d749=\", \"35\u00d762=\", \"16\u00d785=\", \"39\u00d749=\", \"69\u00d727=\"],\n};\n\nfor (const rowIndex of Object.keys(rowUpdates)) {\n  const values = rowUpdates[rowIndex];\n  for (let col = 0; col < values.length; col++) {\n    const cell = table.getCell(Number(rowIndex), col);\n    cell.value = values[col];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and the 25 multiplication problems.\n# The table's rows 1, 5, 10, 15, 20 (1-indexed, Word COM style) are the only\n# rows that carry the \"AxB=\" text cells -- the rows in between are blank\n# spacer rows.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date/day-of-week heading.\n$find = $d.Content.Find\n$find.Text = \"2023-11-10 Friday\"\n$find.Replacement.Text = \"2023-11-11 Saturday\"\n$find.Execute($null, $false, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n\n# 2) Update the multiplication problems table.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1, 1).Range.Text = \"46\u00d711=\"\n$t.Cell(1, 2).Range.Text = \"93\u00d745=\"\n$t.Cell(1, 3).Range.Text = \"27\u00d782=\"\n$t.Cell(1, 4).Range.Text = \"35\u00d775=\"\n$t.Cell(1, 5).Range.Text = \"87\u00d763=\"\n\n$t.Cell(5, 1).Range.Text = \"19\u00d778=\"\n$t.Cell(5, 2).Range.Text = \"26\u00d769=\"\n$t.Cell(5, 3).Range.Text = \"17\u00d784=\"\n$t.Cell(5, 4).Range.Text = \"12\u00d744=\"\n$t.Cell(5, 5).Range.Text = \"71\u00d762=\"\n\n$t.Cell(10, 1).Range.Text = \"45\u00d732=\"\n$t.Cell(10, 2).Range.Text = \"74\u00d713=\"\n$t.Cell(10, 3).Range.Text = \"88\u00d772=\"\n$t.Cell(10, 4).Range.Text = \"29\u00d757=\"\n$t.Cell(10, 5).Range.Text = \"18\u00d752=\"\n\n$t.Cell(15, 1).Range.Text = \"97\u00d750=\"\n$t.Cell(15, 2).Range.Text = \"96\u00d798=\"\n$t.Cell(15, 3).Range.Text = \"25\u00d741=\"\n$t.Cell(15, 4).Range.Text = \"81\u00d744=\"\n$t.Cell(15, 5).Range.Text = \"75\u00d779=\"\n\n$t.Cell(20, 1).Range.Text = \"37\u00d749=\"\n$t.Cell(20, 2).Range.Text = \"35\u00d762=\"\n$t.Cell(20, 3).Range.Text = \"16\u00d785=\"\n$t.Cell(20, 4).Range.Text = \"39\u00d749=\"\n$t.Cell(20, 5).Range.Text = \"69\u00d727=\"\n"}
